$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A23").Value = 21
$ws.Range("A24").Value = 22

$ws.Range("C23").Value = "c=a<b;"
$ws.Range("C24").Value = "c=a>b;"

$ws.Range("B23").Value = "smaller than"
$ws.Range("B24").Value = "larget than"

$ws.Range("B24").Select()
